$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7202
$ws.Range("F7").Value = 180
$ws.Range("F8").Value = 129
$ws.Range("F11").Value = 60
$ws.Range("F12").Value = 216
$ws.Range("F13").Value = 12
$ws.Range("F14").Value = 459
$ws.Range("F16").Value = 1855
$ws.Range("F17").Value = 49
$ws.Range("F18").Value = 35
$ws.Range("F19").Value = 3751
$ws.Range("F21").Value = 250
$ws.Range("F23").Value = 37
$ws.Range("F25").Value = 34
$ws.Range("F26").Value = 2414
$ws.Range("F27").Value = 20
$ws.Range("F28").Value = 295
$ws.Range("F30").Value = 5
$ws.Range("F31").Value = 41
$ws.Range("F33").Value = 16
$ws.Range("F36").Value = 26
$ws.Range("F37").Value = 164
$ws.Range("F38").Value = 1438
$ws.Range("F39").Value = 146

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 8

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7202
$ws.Range("F6").Value = 8
$ws.Range("F8").Value = 180
$ws.Range("F9").Value = 129
$ws.Range("F12").Value = 60
$ws.Range("F13").Value = 216
$ws.Range("F14").Value = 12
$ws.Range("F15").Value = 459
$ws.Range("F17").Value = 1855
$ws.Range("F18").Value = 49
$ws.Range("F19").Value = 35
$ws.Range("F20").Value = 3751
$ws.Range("F22").Value = 250
$ws.Range("F24").Value = 37
$ws.Range("F26").Value = 34
$ws.Range("F27").Value = 2414
$ws.Range("F28").Value = 20
$ws.Range("F29").Value = 295
$ws.Range("F31").Value = 5
$ws.Range("F32").Value = 41
$ws.Range("F34").Value = 16
$ws.Range("F37").Value = 26
$ws.Range("F38").Value = 164
$ws.Range("F39").Value = 1438
$ws.Range("F40").Value = 146
